$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts old D..G -> E..H)
$ws.Columns.Item(4).Insert()

# Remove the old EIRP columns, which are now G (Uplink EIRP) and H (Downlink EIRP)
$ws.Range("G1:H1").EntireColumn.Delete()

# Update header row
$ws.Range("C1").Value = "Distancia Uplink mx (km)"
$ws.Range("D1").Value = "Distancia Downlink mx (km)"
$ws.Range("E1").Value = "Uplink Eb/No (dB)"
$ws.Range("F1").Value = "Downlink Eb/No (dB)"

# Update data rows 2-9 with new values
$distUplink = 7248.583046900511
$distDownlink = 11233.88175375904
$uplinkEbNo = 0.6786769385847222
$downlinkEbNo = 3.885492990715338

for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 3).Value = $distUplink
    $ws.Cells.Item($r, 4).Value = $distDownlink
    $ws.Cells.Item($r, 5).Value = $uplinkEbNo
    $ws.Cells.Item($r, 6).Value = $downlinkEbNo
}
